$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 658.79
$ws.Range("I15").Value = 658.79
$ws.Range("K15").Value = 1976.37
$ws.Range("M15").Value = -1807.37
$ws.Range("H64").Value = 3462.0571
$ws.Range("I64").Value = 3222.68
$ws.Range("J64").Value = 4060.5
$ws.Range("K64").Value = 3222.68
$ws.Range("L64").Value = 4060.5
$ws.Range("M64").Value = -2974.68
$ws.Range("N64").Value = -4556.5
$ws.Range("H67").Value = 3462.0571
$ws.Range("I67").Value = 3222.68
$ws.Range("J67").Value = 4060.5
$ws.Range("K67").Value = 3222.68
$ws.Range("L67").Value = 4060.5
$ws.Range("M67").Value = -2364.68
$ws.Range("N67").Value = -5776.5
$ws.Range("H132").Value = 1956.875
$ws.Range("I132").Value = 1917.5518
$ws.Range("K132").Value = 5752.6554
$ws.Range("M132").Value = -3222.6554
$ws.Range("H137").Value = 3183.7188
$ws.Range("I137").Value = 1774.5
$ws.Range("J137").Value = 6284
$ws.Range("K137").Value = 5323.5
$ws.Range("L137").Value = 18852
$ws.Range("M137").Value = -2773.5
$ws.Range("N137").Value = -23952
$ws.Range("H138").Value = 1675.18
$ws.Range("I138").Value = 760
$ws.Range("J138").Value = 1980.24
$ws.Range("K138").Value = 2280
$ws.Range("L138").Value = 5940.72
$ws.Range("M138").Value = 2860
$ws.Range("N138").Value = -16220.72

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6112.881
$ws.Range("I61").Value = 5797.909
$ws.Range("J61").Value = 7267.778
$ws.Range("K61").Value = 5797.909
$ws.Range("L61").Value = 7267.778
$ws.Range("M61").Value = -5585.909
$ws.Range("N61").Value = -7691.778
$ws.Range("H74").Value = 4367.875
$ws.Range("I74").Value = 3756.45
$ws.Range("K74").Value = 3756.45
$ws.Range("M74").Value = -2882.45
$ws.Range("H77").Value = 4367.875
$ws.Range("I77").Value = 3756.45
$ws.Range("K77").Value = 18782.25
$ws.Range("M77").Value = -14414.25
$ws.Range("H136").Value = 6112.881
$ws.Range("I136").Value = 5797.909
$ws.Range("J136").Value = 7267.778
$ws.Range("K136").Value = 17393.727
$ws.Range("L136").Value = 21803.334
$ws.Range("M136").Value = -14843.727
$ws.Range("N136").Value = -26903.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1430.138
$ws.Range("I94").Value = 1109.1578
$ws.Range("J94").Value = 2040
$ws.Range("K94").Value = 1109.1578
$ws.Range("L94").Value = 2040
$ws.Range("M94").Value = -658.1578
$ws.Range("N94").Value = -2942
$ws.Range("H134").Value = 4735.1143
$ws.Range("I134").Value = 4885.303
$ws.Range("J134").Value = 2257
$ws.Range("K134").Value = 14655.909
$ws.Range("L134").Value = 6771
$ws.Range("M134").Value = -12120.909
$ws.Range("N134").Value = -11841

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2600650.8
$ws.Range("I58").Value = 5683871
$ws.Range("K58").Value = 5683871
$ws.Range("M58").Value = -5683668
$ws.Range("H136").Value = 2600650.8
$ws.Range("I136").Value = 5683871
$ws.Range("K136").Value = 17051613
$ws.Range("M136").Value = -17049063

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 125001864
$ws.Range("I22").Value = 333333630
$ws.Range("J22").Value = 2800
$ws.Range("K22").Value = 1000000890
$ws.Range("L22").Value = 8400
$ws.Range("M22").Value = -1000000721
$ws.Range("N22").Value = -8738
$ws.Range("H23").Value = 1206.6666
$ws.Range("I23").Value = 10001
$ws.Range("J23").Value = 107.375
$ws.Range("K23").Value = 30003
$ws.Range("L23").Value = 322.125
$ws.Range("M23").Value = -29768
$ws.Range("N23").Value = -792.125
$ws.Range("H27").Value = 125001864
$ws.Range("I27").Value = 333333630
$ws.Range("J27").Value = 2800
$ws.Range("K27").Value = 1000000890
$ws.Range("L27").Value = 8400
$ws.Range("M27").Value = -1000000788
$ws.Range("N27").Value = -8604
$ws.Range("H29").Value = 207.5
$ws.Range("I29").Value = 80
$ws.Range("J29").Value = 250
$ws.Range("K29").Value = 240
$ws.Range("L29").Value = 750
$ws.Range("M29").Value = 37
$ws.Range("N29").Value = -1304
$ws.Range("H36").Value = 1000
$ws.Range("J36").Value = 2000
$ws.Range("L36").Value = 6000
$ws.Range("N36").Value = -6338
$ws.Range("H46").Value = 3029.1667
$ws.Range("J46").Value = 3191.1765
$ws.Range("L46").Value = 9573.529500000001
$ws.Range("N46").Value = -9755.529500000001
$ws.Range("H58").Value = 2298.8
$ws.Range("I58").Value = 392.5
$ws.Range("J58").Value = 2775.375
$ws.Range("K58").Value = 1177.5
$ws.Range("L58").Value = 8326.125
$ws.Range("M58").Value = -1049.5
$ws.Range("N58").Value = -8582.125
$ws.Range("H86").Value = 3300.6667
$ws.Range("I86").Value = 3300.6667
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 9902.000100000001
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -8716.000100000001
$ws.Range("N86").Value = $null
$ws.Range("H89").Value = 3300.6667
$ws.Range("I89").Value = 3300.6667
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 29706.0003
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -23778.0003
$ws.Range("N89").Value = -35358
$ws.Range("H132").Value = 2053.4443
$ws.Range("I132").Value = 2906.8
$ws.Range("J132").Value = 1725.2307
$ws.Range("K132").Value = 26161.2
$ws.Range("L132").Value = 15527.0763
$ws.Range("M132").Value = -23631.2
$ws.Range("N132").Value = -20587.0763

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 78.95238000000001
$ws.Range("I2").Value = 39.125
$ws.Range("J2").Value = 206.4
$ws.Range("K2").Value = 39.125
$ws.Range("L2").Value = 206.4
$ws.Range("M2").Value = 73.875
$ws.Range("N2").Value = -432.4
$ws.Range("H136").Value = 9635.583000000001
$ws.Range("J136").Value = 9635.583000000001
$ws.Range("L136").Value = 28906.749
$ws.Range("N136").Value = -34006.749

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2283.5
$ws.Range("I16").Value = 2283.5
$ws.Range("K16").Value = 2283.5
$ws.Range("M16").Value = -2113.5
$ws.Range("H93").Value = 1041.2
$ws.Range("I93").Value = 1041.2
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1041.2
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 206.8
$ws.Range("N93").Value = $null
$ws.Range("H132").Value = 3654.9473
$ws.Range("I132").Value = 2802.0908
$ws.Range("J132").Value = 4827.625
$ws.Range("K132").Value = 8406.2724
$ws.Range("L132").Value = 14482.875
$ws.Range("M132").Value = -5876.2724
$ws.Range("N132").Value = -19542.875
$ws.Range("H133").Value = 34332.285
$ws.Range("J133").Value = 34332.285
$ws.Range("L133").Value = 34332.285
$ws.Range("N133").Value = -39392.285

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1782.5
$ws.Range("I81").Value = 1090
$ws.Range("J81").Value = 2475
$ws.Range("K81").Value = 2180
$ws.Range("L81").Value = 4950
$ws.Range("M81").Value = -1119
$ws.Range("N81").Value = -7072
$ws.Range("H84").Value = 1782.5
$ws.Range("I84").Value = 1090
$ws.Range("J84").Value = 2475
$ws.Range("K84").Value = 10900
$ws.Range("L84").Value = 24750
$ws.Range("M84").Value = -5596
$ws.Range("N84").Value = -35358
